$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all the G column values (rows 5-18, 20-26) to 51, leaving G3, G4, G19 untouched
$rows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,20,21,22,23,24,25,26)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = 51
}

# Update the selected cell/active cell to G24
$ws.Range("G24").Select()
